# Fill in the "Amount" column (D) for 2025-01-10 .. 2025-01-22 (rows 17-29)
# with the per-day hormone amounts, matching the style already used by the
# other filled-in Amount cells (D2:D16), i.e. style index 1 (wrap text /
# vertical-top alignment).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$amounts = @{
    17 = "5`n10`n3"
    18 = "5`n10`n2"
    19 = "6`n12`n1"
    20 = "6`n12`n1"
    21 = "6`n14`n1"
    22 = "6`n16`n1"
    23 = "6`n14`n1"
    24 = "6`n12`n1"
    25 = "6`n12`n1"
    26 = "6`n10`n1"
    27 = "6`n10`n1"
    28 = "6`n6`n1"
    29 = "6`n6`n1"
}

# Grab the format (wrap text + vertical-top alignment) already used by the
# existing Amount entries, e.g. D2, so the newly written cells share the
# same style instead of creating a brand new one.
$ws.Range("D2").Copy()

foreach ($row in 17..29) {
    $cell = $ws.Range("D$row")
    $cell.PasteSpecial(-4122)
    $cell.Value = $amounts[$row]
}
